$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsEHP = $wb.Worksheets.Item("EHPpUC")

# --- Calculations sheet -------------------------------------------------

# Insert a new (blank) row above the old row 10 ("lb per metric ton" block)
# so the lower section shifts down by one row; this is where the new
# commentary text will live underneath the corrected HHV figure.
$wsCalc.Rows.Item(9).Insert()

# Bug fix: correct Hydrogen HHV from 23874 BTU/lb to 60920 BTU/lb.
$wsCalc.Range("A8").Value() = 60920

# New reviewer commentary notes next to the "lb per metric ton" /
# "BTU H2 / MW (annual production)" rows (now shifted to rows 12-14).
$wsCalc.Range("C12").Value() = "The study they are citing uses 3 scenarios of Fuel Cell Electric Vehicle adoption. Then it calculates the amount of hydrogen needed to suppor those vehicles. Then it calculates the electrolyzer capacity needed to supply that hydrogen. "
$wsCalc.Range("C13").Value() = "So, I think it's fair, using the EPS assumptions of 24/7/365 operation, that smallest electrolyzer you would need to produce 1.39e10 annual Btu would be 1 MW."
$wsCalc.Range("C14").Value() = "No reason to think this would be different for Texas."

$wsCalc.Range("C12:C14").Font.ThemeColor() = 9

# --- Selections (cosmetic, matches the saved cursor positions) ----------

$wsCalc.Range("J6").Select()
$wsEHP.Range("B2").Select()
$wsAbout.Range("B6").Select()
$wsAbout.Activate()

Write-Output "done"
